# Updates "horarios" workbook with a fresh scrape timestamp (06:03:38)
# for line 141, adding newly-arrived rows to each of the three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:03:38"
$ws1.Range("A3").Value = "Total filas: 32"

$sheet1Rows = @(
    @("15","06:03:38","06:11","215A_EL PATO",8,"LP1912"),
    @("16","06:03:38","06:14","225_HARAS DEL SUR",11,"LP1912"),
    @("17","06:03:38","06:21","26_HERNANDEZ",18,"LP1912"),
    @("18","06:03:38","06:27","23_HERNANDEZ",24,"LP1912"),
    @("19","06:03:38","06:29","86_EST CHICA-ESC AGRARIA",26,"LP1912"),
    @("20","06:03:38","06:31","16_SANTA ANA",28,"LP1912"),
    @("21","06:03:38","06:44","225_C ROCA-H SUR",41,"LP1912"),
    @("22","06:03:38","06:46","215C_EL PATO",43,"LP1912"),
    @("23","06:03:38","06:59","14_ABASTO",56,"LP1912"),
    @("24","06:03:38","07:04","23_HERNANDEZ",61,"LP1912"),
    @("25","06:03:38","07:05","15_ABASTO",62,"LP1912"),
    @("26","06:03:38","07:07","225_GOMEZ",64,"LP1912"),
    @("27","06:03:38","07:11","215A_EL PATO",68,"LP1912"),
    @("28","06:03:38","07:15","11_ETCHEVERRY",72,"LP1912"),
    @("29","06:03:38","07:21","26_HERNANDEZ",78,"LP1912"),
    @("30","06:03:38","07:23","10_OLMOS",80,"LP1912"),
    @("31","06:03:38","07:31","11_ETCHEVERRY",88,"LP1912"),
    @("32","06:03:38","07:31","16_SANTA ANA",88,"LP1912"),
    @("33","06:03:38","07:32","84_COLONIA URQUIZA-ESC 49",89,"LP1912"),
    @("34","06:03:38","07:36","27_EL RETIRO",93,"LP1912"),
    @("35","06:03:38","07:39","10_OLMOS",96,"LP1912"),
    @("36","06:03:38","07:47","14_ABASTO",104,"LP1912"),
    @("37","06:03:38","07:51","215D_EL PATO",108,"LP1912")
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:03:38"
$ws2.Range("A3").Value = "Total filas: 6"

$sheet2Rows = @(
    @("8","06:03:38","06:11","215A_EL PATO",8,"LP1912"),
    @("9","06:03:38","06:46","215C_EL PATO",43,"LP1912"),
    @("10","06:03:38","07:11","215A_EL PATO",68,"LP1912"),
    @("11","06:03:38","07:51","215D_EL PATO",108,"LP1912")
)

foreach ($row in $sheet2Rows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 06:03:38"
$ws3.Range("A3").Value = "Total filas: 7"

$sheet3Rows = @(
    @("8","06:03:38","06:08","215A_LA PLATA",5,"L6173"),
    @("9","06:03:38","06:32","215C_LA PLATA",29,"L6203"),
    @("11","06:03:38","07:00","215B_LP-P MOR-1 Y 57",57,"L6173"),
    @("12","06:03:38","07:35","215A_LA PLATA",92,"L6173")
)

foreach ($row in $sheet3Rows) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
}
